$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the calibration values for L1 (row 4)
$ws.Range("F4").Value = 1950
$ws.Range("G4").Value = 2050

# Update the window view: scroll back to A1 (remove topLeftCell offset) and change selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H12").Select()
